# Update template files for import trips
# - "Facility" header (D1) becomes "Facility*" (now a required field)
# - Active cell selection moves to D2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Facility*"

$ws.Range("D2").Select()
